# Add a new timesheet entry (row 79) to the "hours" worksheet, mirroring
# the existing pattern: date (col A, same date format as the rows above),
# hours worked (col B), and the running-total formula (col C) that
# continues on from C78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from the last existing row (A78) onto the new
# row's date cell so it reuses the same cell style instead of Excel
# minting a brand-new (duplicate) number format.
$ws.Range("A78").Copy()
$ws.Range("A79").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New entry: 2024-12-23 (serial 45649), 1 hour worked.
$ws.Range("A79").Value = 45649
$ws.Range("B79").Value = 1

# Continue the running-total formula pattern from C78.
$ws.Range("C79").Formula = "=C78+B79"

# Leave the selection on the newly-entered cell, as Excel would after
# the user finishes typing the new row.
$ws.Range("C79").Select()
